$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @("ECs", "Bmp4", "Bmpr2", "ECs", 2, 0.6666666666666666, 4.504599666666667, 13.513799, 0.2870666703033801, 0.2870666703033801, 3, 1, 33.66677566666667, 101.000327, 0.26047027611676, 0.26047027611676, 151.6553464458081, 1364.898118012273, 0.07477233487784032, 0.07477233487784034)
  3 = @("ECs", "Bmp4", "Bmpr2", "FAPs", 2, 0.6666666666666666, 4.504599666666667, 13.513799, 0.2870666703033801, 0.2870666703033801, 3, 1, 43.683024, 131.049072, 0.3379631431162116, 0.3379631431162116, 196.774535349392, 1770.970818144528, 0.09701795417963556, 0.09701795417963557)
  4 = @("ECs", "Bmp4", "Bmpr2", "M2", 2, 0.6666666666666666, 4.504599666666667, 13.513799, 0.2870666703033801, 0.2870666703033801, 3, 1, 17.657557, 52.97267100000001, 0.1366115006935798, 0.1366115006935798, 79.54022537634766, 715.8620283871292, 0.03921660862925386, 0.03921660862925386)
  5 = @("ECs", "Bmp4", "Bmpr2", "sCs", 2, 0.6666666666666666, 4.504599666666667, 13.513799, 0.2870666703033801, 0.2870666703033801, 3, 1, 34.24645366666667, 102.739361, 0.2649550800734486, 0.2649550800734485, 154.2665637713821, 1388.399073942439, 0.07605977261665033, 0.07605977261665033)
  6 = @("FAPs", "Bmp4", "Bmpr2", "ECs", 3, 1, 7.480073, 22.440219, 0.476686011772755, 0.476686011772755, 3, 1, 33.66677566666667, 101.000327, 0.26047027611676, 0.26047027611676, 251.8299396612903, 2266.469456951613, 0.1241625371074466, 0.1241625371074466)
  7 = @("FAPs", "Bmp4", "Bmpr2", "FAPs", 3, 1, 7.480073, 22.440219, 0.476686011772755, 0.476686011772755, 3, 1, 43.683024, 131.049072, 0.3379631431162116, 0.3379631431162116, 326.752208380752, 2940.769875426768, 0.1611023028182517, 0.1611023028182517)
  8 = @("FAPs", "Bmp4", "Bmpr2", "M2", 3, 1, 7.480073, 22.440219, 0.476686011772755, 0.476686011772755, 3, 1, 17.657557, 52.97267100000001, 0.1366115006935798, 0.1366115006935798, 132.079815361661, 1188.718338254949, 0.06512079142791354, 0.06512079142791354)
  9 = @("FAPs", "Bmp4", "Bmpr2", "sCs", 3, 1, 7.480073, 22.440219, 0.476686011772755, 0.476686011772755, 3, 1, 34.24645366666667, 102.739361, 0.2649550800734486, 0.2649550800734485, 256.1659734177844, 2305.493760760059, 0.1263003804191432, 0.1263003804191432)
  10 = @("sCs", "Bmp4", "Bmpr2", "ECs", 3, 1, 3.707151333333333, 11.121454, 0.2362473179238649, 0.2362473179238649, 3, 1, 33.66677566666667, 101.000327, 0.26047027611676, 0.26047027611676, 124.8078323017176, 1123.270490715458, 0.06153540413147308, 0.06153540413147309)
  11 = @("sCs", "Bmp4", "Bmpr2", "FAPs", 3, 1, 3.707151333333333, 11.121454, 0.2362473179238649, 0.2362473179238649, 3, 1, 43.683024, 131.049072, 0.3379631431162116, 0.3379631431162116, 161.939580665632, 1457.456225990688, 0.07984288611832428, 0.0798428861183243)
  12 = @("sCs", "Bmp4", "Bmpr2", "M2", 3, 1, 3.707151333333333, 11.121454, 0.2362473179238649, 0.2362473179238649, 3, 1, 17.657557, 52.97267100000001, 0.1366115006935798, 0.1366115006935798, 65.45923597595934, 589.133123783634, 0.03227410063641244, 0.03227410063641244)
  13 = @("sCs", "Bmp4", "Bmpr2", "sCs", 3, 1, 3.707151333333333, 11.121454, 0.2362473179238649, 0.2362473179238649, 3, 1, 34.24645366666667, 102.739361, 0.2649550800734486, 0.2649550800734485, 126.9567863723216, 1142.611077350894, 0.06259492703765508, 0.06259492703765508)
}

foreach ($r in $data.Keys) {
  $row = $data[$r]
  for ($i = 0; $i -lt $row.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item([int]$r, $col).Value = $row[$i]
  }
}
